# 1.1.6_assess_by_word.xlsx -- "File updates on program run"
#
# The backend job re-ran: it (a) logged a brand-new word ("ਲੇਪਨ") into the
# Words sheet (selected + analysis already kicked off), (b) re-stamped the
# "selected_at" timestamp for the three existing ਲੇਪਨ occurrences already
# queued in Progress (rows 5-7), (c) finished analysis for two of the
# occurrences that were already sitting in Progress (rows 2, 3, 5 & 6), and
# (d) appended a fresh trio of Progress rows (8-10) for the new word, two of
# which finish analysis immediately and one which is left "not started".
#
# The timestamps below are the literal run-time values the program wrote to
# disk for this pass (they are data captured from that run, not something a
# formula can re-derive).

$wb = $excel.ActiveWorkbook
$wsWords    = $wb.Worksheets.Item("Words")
$wsProgress = $wb.Worksheets.Item("Progress")

$dateFmt = "YYYY-MM-DD HH:MM:SS"

$word = "ਲੇਪਨ"

# Moments captured during this run.
$tListed      = 45912.28803851852   # word logged / newly-queued rows selected
$tStarted     = 45912.28809836806   # analysis kicked off for the new word
$tReselected  = 45912.27207552084   # rows 5-7 re-selected for analysis
$tDoneFirst   = 45912.28862107639   # 1st occurrence analysis completed
$tDoneSecond  = 45912.28908845241   # 2nd occurrence analysis completed

# ---------------------------------------------------------------------
# Words sheet: one new row for the newly-logged word.
# ---------------------------------------------------------------------
$wRow = 2
$wsWords.Cells.Item($wRow, 1).Value  = $word          # word
$wsWords.Cells.Item($wRow, 2).Value  = $word          # word_key_norm
$wsWords.Cells.Item($wRow, 3).Value  = $true          # listed_by_user
$wsWords.Cells.Item($wRow, 4).Value  = $tListed       # listed_at
$wsWords.Cells.Item($wRow, 4).NumberFormat = $dateFmt
$wsWords.Cells.Item($wRow, 5).Value  = $true          # selected_for_analysis
$wsWords.Cells.Item($wRow, 6).Value  = $tListed       # selected_at
$wsWords.Cells.Item($wRow, 6).NumberFormat = $dateFmt
$wsWords.Cells.Item($wRow, 7).Value  = $true          # analysis_started
$wsWords.Cells.Item($wRow, 8).Value  = $tStarted      # analysis_started_at
$wsWords.Cells.Item($wRow, 8).NumberFormat = $dateFmt
$wsWords.Cells.Item($wRow, 9).Value  = $false         # analysis_completed
$wsWords.Cells.Item($wRow, 10).Value = ""             # analysis_completed_at (blank)
$wsWords.Cells.Item($wRow, 11).Value = 0              # sequence_index
$wsWords.Cells.Item($wRow, 12).Value = ""             # notes (blank)

# ---------------------------------------------------------------------
# Progress sheet: re-selection stamp for the three rows already queued
# for this word (rows 5, 6, 7).
# ---------------------------------------------------------------------
5..7 | ForEach-Object {
    $r = $_
    $wsProgress.Cells.Item($r, 7).Value = $tReselected   # selected_at
    $wsProgress.Cells.Item($r, 7).NumberFormat = $dateFmt
}

# Normalise page_number to a real number on the three re-selected rows
# (they had previously been written as text).
$wsProgress.Cells.Item(5, 5).Value = 1018
$wsProgress.Cells.Item(6, 5).Value = 1103
$wsProgress.Cells.Item(7, 5).Value = 1243

# First occurrence of the word finishes analysis (existing row 2 and the
# matching re-selected row 5).
2, 5 | ForEach-Object {
    $r = $_
    $wsProgress.Cells.Item($r, 8).Value = "completed"
    $wsProgress.Cells.Item($r, 9).Value = $tDoneFirst
    $wsProgress.Cells.Item($r, 9).NumberFormat = $dateFmt
}

# Second occurrence of the word finishes analysis (existing row 3 and the
# matching re-selected row 6).
3, 6 | ForEach-Object {
    $r = $_
    $wsProgress.Cells.Item($r, 8).Value = "completed"
    $wsProgress.Cells.Item($r, 9).Value = $tDoneSecond
    $wsProgress.Cells.Item($r, 9).NumberFormat = $dateFmt
}

# Row 7 (third occurrence) stays "not started" -- untouched otherwise.

# ---------------------------------------------------------------------
# Progress sheet: three brand-new rows (8, 9, 10) for the newly-logged
# word, mirroring the verses already tracked for it.
# ---------------------------------------------------------------------
$verses = @(
    @{ Row = 8;  Verse = "ਚੰਦਨ ਅਗਰ ਕਪੂਰ ਲੇਪਨ ਤਿਸੁ ਸੰਗੇ ਨਹੀ ਪ੍ਰੀਤਿ ॥"; Page = 1018; Status = "completed";   Done = $tDoneFirst },
    @{ Row = 9;  Verse = "ਜਟਾ ਭਸਮ ਲੇਪਨ ਕੀਆ ਕਹਾ ਗੁਫਾ ਮਹਿ ਬਾਸੁ ॥";       Page = 1103; Status = "completed";   Done = $tDoneSecond },
    @{ Row = 10; Verse = "ਬਾਹਰਿ ਭਸਮ ਲੇਪਨ ਕਰੇ ਅੰਤਰਿ ਗੁਬਾਰੀ ॥";           Page = 1243; Status = "not started"; Done = $null }
)

foreach ($entry in $verses) {
    $r = $entry.Row
    $wsProgress.Cells.Item($r, 1).Value  = $word            # word
    $wsProgress.Cells.Item($r, 2).Value  = $word            # word_key_norm
    $wsProgress.Cells.Item($r, 3).Value  = ""               # word_index (blank)
    $wsProgress.Cells.Item($r, 4).Value  = $entry.Verse     # verse
    $wsProgress.Cells.Item($r, 5).Value  = $entry.Page      # page_number
    $wsProgress.Cells.Item($r, 6).Value  = $true            # selected_for_analysis
    $wsProgress.Cells.Item($r, 7).Value  = $tListed         # selected_at
    $wsProgress.Cells.Item($r, 7).NumberFormat = $dateFmt
    $wsProgress.Cells.Item($r, 8).Value  = $entry.Status    # status
    if ($entry.Done -ne $null) {
        $wsProgress.Cells.Item($r, 9).Value = $entry.Done    # completed_at
        $wsProgress.Cells.Item($r, 9).NumberFormat = $dateFmt
    } else {
        $wsProgress.Cells.Item($r, 9).Value = ""              # completed_at (blank)
    }
    $wsProgress.Cells.Item($r, 10).Value = ""               # reanalyzed_count (blank)
    $wsProgress.Cells.Item($r, 11).Value = ""               # last_reanalyzed_at (blank)
}
